$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C to Text format so date-like strings are not reinterpreted as dates
$ws.Range("C2:C70").NumberFormat = "@"

$data = @(
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2019", 87.15000000000001)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2019", 87.86)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2019", 88.09999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2019", 88.92)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2020", 87.63)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2020", 86.40000000000001)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2020", 85.11)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2020", 85.81999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2021", 85.09)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2021", 85.77)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2021", 87.36)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2021", 88.84999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2022", 88.86)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2022", 90.7)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2022", 91.3)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2022", 92.06)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2023", 91.20999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2023", 91.95999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2023", 92.31)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2023", 92.59)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2024", 92.08)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2024", 93.09999999999999)
    ,("Brasil", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2024", 93.64)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2019", 84.55)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2019", 85.23999999999999)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2019", 85.39)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2019", 86.23999999999999)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2020", 84.22)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2020", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2020", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2020", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2021", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2021", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2021", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2021", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2022", $null)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2022", 87.3)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2022", 88.03)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2022", 89.14)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2023", 87.76000000000001)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2023", 88.67)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2023", 89.15000000000001)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2023", 89.56)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2024", 88.89)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2024", 90.56999999999999)
    ,("Nordeste", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2024", 91.33)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2019", 84.48)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2019", 84.65000000000001)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2019", 85.22)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2019", 85.06999999999999)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2020", 84.20999999999999)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2020", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2020", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2020", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2021", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2021", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2021", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2021", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2022", $null)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2022", 87.25)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2022", 87.87)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2022", 88.05)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2023", 88.17)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2023", 89.76000000000001)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2023", 90.2)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/10/2023", 88.7)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/01/2024", 89.93000000000001)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/04/2024", 90.95999999999999)
    ,("Sergipe", "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência", "01/07/2024", 91.61)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 4).Value = ""
    }
    $r++
}

# Reset column C formatting back to default (General) now that text is locked in as a string type
$ws.Range("C2:C70").ClearFormats()

# Re-apply header style (ClearFormats above only touched C2:C70, but do the header border+alignment now
# so it is independent of ordering)
$ws.Range("A1:D1").Borders.LineStyle = 1
$ws.Range("A1:D1").Borders.Weight = 2
$ws.Range("A1:D1").VerticalAlignment = -4160

[void]$ws.Range("A1").Select()
